$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Slpi"
$ws.Cells.Item(2,3).Value = "Plscr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 27.47303533333333
$ws.Cells.Item(2,8).Value = 82.419106
$ws.Cells.Item(2,9).Value = 0.9778496847444623
$ws.Cells.Item(2,10).Value = 0.9778496847444623
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 10.38882
$ws.Cells.Item(2,14).Value = 31.16646
$ws.Cells.Item(2,15).Value = 0.6821017765484108
$ws.Cells.Item(2,16).Value = 0.6821017765484106
$ws.Cells.Item(2,17).Value = 285.41241893164
$ws.Cells.Item(2,18).Value = 2568.71177038476
$ws.Cells.Item(2,19).Value = 0.6669930071615011
$ws.Cells.Item(2,20).Value = 0.666993007161501

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Slpi"
$ws.Cells.Item(3,3).Value = "Plscr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 27.47303533333333
$ws.Cells.Item(3,8).Value = 82.419106
$ws.Cells.Item(3,9).Value = 0.9778496847444623
$ws.Cells.Item(3,10).Value = 0.9778496847444623
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.786435
$ws.Cells.Item(3,14).Value = 11.359305
$ws.Cells.Item(3,15).Value = 0.248607064159845
$ws.Cells.Item(3,16).Value = 0.248607064159845
$ws.Cells.Item(3,17).Value = 104.02486254237
$ws.Cells.Item(3,18).Value = 936.2237628813299
$ws.Cells.Item(3,19).Value = 0.2431003393139508
$ws.Cells.Item(3,20).Value = 0.2431003393139508

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Slpi"
$ws.Cells.Item(4,3).Value = "Plscr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 27.47303533333333
$ws.Cells.Item(4,8).Value = 82.419106
$ws.Cells.Item(4,9).Value = 0.9778496847444623
$ws.Cells.Item(4,10).Value = 0.9778496847444623
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.055346
$ws.Cells.Item(4,14).Value = 3.166038
$ws.Cells.Item(4,15).Value = 0.06929115929174429
$ws.Cells.Item(4,16).Value = 0.0692911592917443
$ws.Cells.Item(4,17).Value = 28.993557946892
$ws.Cells.Item(4,18).Value = 260.942021522028
$ws.Cells.Item(4,19).Value = 0.06775633826901047
$ws.Cells.Item(4,20).Value = 0.06775633826901048

# Row 5
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Slpi"
$ws.Cells.Item(5,3).Value = "Plscr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.622321
$ws.Cells.Item(5,8).Value = 1.866963
$ws.Cells.Item(5,9).Value = 0.02215031525553766
$ws.Cells.Item(5,10).Value = 0.02215031525553766
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 10.38882
$ws.Cells.Item(5,14).Value = 31.16646
$ws.Cells.Item(5,15).Value = 0.6821017765484108
$ws.Cells.Item(5,16).Value = 0.6821017765484106
$ws.Cells.Item(5,17).Value = 6.46518085122
$ws.Cells.Item(5,18).Value = 58.18662766098001
$ws.Cells.Item(5,19).Value = 0.01510876938690961
$ws.Cells.Item(5,20).Value = 0.0151087693869096

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Slpi"
$ws.Cells.Item(6,3).Value = "Plscr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.622321
$ws.Cells.Item(6,8).Value = 1.866963
$ws.Cells.Item(6,9).Value = 0.02215031525553766
$ws.Cells.Item(6,10).Value = 0.02215031525553766
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.786435
$ws.Cells.Item(6,14).Value = 11.359305
$ws.Cells.Item(6,15).Value = 0.248607064159845
$ws.Cells.Item(6,16).Value = 0.248607064159845
$ws.Cells.Item(6,17).Value = 2.356378015635
$ws.Cells.Item(6,18).Value = 21.207402140715
$ws.Cells.Item(6,19).Value = 0.005506724845894246
$ws.Cells.Item(6,20).Value = 0.005506724845894246

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Slpi"
$ws.Cells.Item(7,3).Value = "Plscr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.622321
$ws.Cells.Item(7,8).Value = 1.866963
$ws.Cells.Item(7,9).Value = 0.02215031525553766
$ws.Cells.Item(7,10).Value = 0.02215031525553766
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.055346
$ws.Cells.Item(7,14).Value = 3.166038
$ws.Cells.Item(7,15).Value = 0.06929115929174429
$ws.Cells.Item(7,16).Value = 0.0692911592917443
$ws.Cells.Item(7,17).Value = 0.656763978066
$ws.Cells.Item(7,18).Value = 5.910875802594
$ws.Cells.Item(7,19).Value = 0.001534821022733814
$ws.Cells.Item(7,20).Value = 0.001534821022733814
